$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Header for new column F
$ws.Range("F1").Value = "time_taken"
# Copy formatting (font, border, alignment) from E1's header style
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# time_taken values for rows 2-17, written as text
$times = @(
    "2021-10-05 10:51:12.723712",
    "2021-10-05 10:51:12.723722",
    "2021-10-05 10:51:12.723725",
    "2021-10-05 10:51:12.723728",
    "2021-10-05 10:51:12.723731",
    "2021-10-05 10:51:12.723733",
    "2021-10-05 10:51:12.723736",
    "2021-10-05 10:51:12.723739",
    "2021-10-05 10:51:12.723742",
    "2021-10-05 10:51:12.723744",
    "2021-10-05 10:51:12.723747",
    "2021-10-05 10:51:12.723749",
    "2021-10-05 10:51:12.723752",
    "2021-10-05 10:51:12.723754",
    "2021-10-05 10:51:12.723757",
    "2021-10-05 10:51:12.723759"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $times[$i]
}
